{"js": "// Apply the text replacements described by the diff:\n//  1. Title (Heading1) + repeated bold title near the end\n//  2. \"What we like\" bullet: winning potential bullet\n//  3. \"What we like\" bullet: free spins / retrigger bullet\n//  4. \"What we don't like\" bullet: theme bullet\n//  5. \"What we don't like\" bullet: backdrop visuals bullet\n//  7. Closing italic summary paragraph\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"Play Hot Blox for Free: Unique Slot Game with Flashy Graphics\",\n    replace: \"Play Hot Blox Free - Innovative Gameplay with Excellent Winning Potential\",\n  },\n  {\n    find: \"Exciting winning potential and multipliers\",\n    replace: \"Excellent winning potential\",\n  },\n  {\n    find: \"Free spins feature with retrigger option\",\n    replace: \"Compatible with desktop and mobile devices\",\n  },\n  {\n    find: \"Uninspired game theme\",\n    replace: \"Lack of clear theme\",\n  },\n  {\n    find: \"Limited backdrop visuals\",\n    replace: \"Limited free spins feature\",\n  },\n  {\n    find: \"Ready to play Hot Blox for free? Read our review to discover this innovative and exciting slot game with flashy graphics and a free spins feature.\",\n    replace: \"Read our review of Hot Blox, a unique slot game with innovative gameplay and free spins feature.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff:\n#  1. Title (Heading1) + repeated bold title near the end\n#  2. \"What we like\" bullet: winning potential bullet\n#  3. \"What we like\" bullet: free spins / retrigger bullet\n#  4. \"What we don't like\" bullet: theme bullet\n#  5. \"What we don't like\" bullet: backdrop visuals bullet\n#  6. Closing italic summary paragraph\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n  @{ Find = \"Play Hot Blox for Free: Unique Slot Game with Flashy Graphics\"; Replace = \"Play Hot Blox Free - Innovative Gameplay with Excellent Winning Potential\" },\n  @{ Find = \"Exciting winning potential and multipliers\"; Replace = \"Excellent winning potential\" },\n  @{ Find = \"Free spins feature with retrigger option\"; Replace = \"Compatible with desktop and mobile devices\" },\n  @{ Find = \"Uninspired game theme\"; Replace = \"Lack of clear theme\" },\n  @{ Find = \"Limited backdrop visuals\"; Replace = \"Limited free spins feature\" },\n  @{ Find = \"Ready to play Hot Blox for free? Read our review to discover this innovative and exciting slot game with flashy graphics and a free spins feature.\"; Replace = \"Read our review of Hot Blox, a unique slot game with innovative gameplay and free spins feature.\" }\n)\n\nforeach ($r in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $r.Replace, $wdReplaceAll)\n}\n\nWrite-Output \"done\"\n"}
